# PowerShell Excel COM-interop script
# Updates recalculated crafting-profit figures (currentAveragePrice,
# NQ/HQ price + profit columns H:N) for the affected leve rows across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets, matching the
# latest market-board pull from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1989.5807
$ws.Range("I17").Value = 1571.4286
$ws.Range("J17").Value = 2111.5417
$ws.Range("K17").Value = 4714.2858
$ws.Range("L17").Value = 6334.625100000001
$ws.Range("M17").Value = -4546.2858
$ws.Range("N17").Value = -6670.625100000001

$ws.Range("H62").Value = 10950.8
$ws.Range("I62").Value = 10251.333
$ws.Range("J62").Value = 12000
$ws.Range("K62").Value = 10251.333
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = -9627.333000000001
$ws.Range("N62").Value = -13248

$ws.Range("H65").Value = 10950.8
$ws.Range("I65").Value = 10251.333
$ws.Range("J65").Value = 12000
$ws.Range("K65").Value = 51256.665
$ws.Range("L65").Value = 60000
$ws.Range("M65").Value = -48136.665
$ws.Range("N65").Value = -66240

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws.Range("H45").Value = 3742.2
$ws.Range("I45").Value = 2737.5
$ws.Range("K45").Value = 2737.5
$ws.Range("M45").Value = -2360.5

$ws.Range("H63").Value = 3421.6428
$ws.Range("I63").Value = 2543.4285
$ws.Range("J63").Value = 4299.857
$ws.Range("K63").Value = 2543.4285
$ws.Range("L63").Value = 4299.857
$ws.Range("M63").Value = -1857.4285
$ws.Range("N63").Value = -5671.857

$ws.Range("H66").Value = 3421.6428
$ws.Range("I66").Value = 2543.4285
$ws.Range("J66").Value = 4299.857
$ws.Range("K66").Value = 12717.1425
$ws.Range("L66").Value = 21499.285
$ws.Range("M66").Value = -9285.1425
$ws.Range("N66").Value = -28363.285

$ws.Range("H102").Value = 8333.333000000001
$ws.Range("I102").Value = 7500
$ws.Range("K102").Value = 7500
$ws.Range("M102").Value = -5878

$ws.Range("H110").Value = 4256.625
$ws.Range("I110").Value = 3962
$ws.Range("K110").Value = 3962
$ws.Range("M110").Value = -1917

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 368.4
$ws.Range("I94").Value = 368.4
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 368.4
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = 82.60000000000002
$ws.Range("N94").ClearContents()

$ws.Range("H105").Value = 1253.0834
$ws.Range("I105").Value = 1233.9
$ws.Range("J105").Value = 1349
$ws.Range("K105").Value = 1233.9
$ws.Range("L105").Value = 1349
$ws.Range("M105").Value = 513.0999999999999
$ws.Range("N105").Value = -4843

$ws.Range("H107").Value = 3301.276
$ws.Range("I107").Value = 1207.7778
$ws.Range("K107").Value = 1207.7778
$ws.Range("M107").Value = 712.2221999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 193.5
$ws.Range("I13").Value = 202
$ws.Range("J13").Value = 163.75
$ws.Range("K13").Value = 202
$ws.Range("L13").Value = 163.75
$ws.Range("M13").Value = -63
$ws.Range("N13").Value = -441.75

$ws.Range("H80").Value = 108449.25
$ws.Range("J80").Value = 108449.25
$ws.Range("L80").Value = 108449.25
$ws.Range("N80").Value = -110695.25

$ws.Range("H83").Value = 108449.25
$ws.Range("J83").Value = 108449.25
$ws.Range("L83").Value = 325347.75
$ws.Range("N83").Value = -336579.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 225
$ws.Range("I57").Value = 225
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 675
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -116
$ws.Range("N57").ClearContents()

$ws.Range("H69").Value = 2305.5
$ws.Range("I69").Value = 1611
$ws.Range("K69").Value = 4833
$ws.Range("M69").Value = -4022

$ws.Range("H72").Value = 2305.5
$ws.Range("I72").Value = 1611
$ws.Range("K72").Value = 14499
$ws.Range("M72").Value = -10443

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H139").Value = 4897.222
$ws.Range("I139").Value = 4477.857
$ws.Range("K139").Value = 13433.571
$ws.Range("M139").Value = -8293.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 7501903
$ws.Range("I7").Value = 14286286
$ws.Range("J7").Value = 2225161
$ws.Range("K7").Value = 14286286
$ws.Range("L7").Value = 2225161
$ws.Range("M7").Value = -14286174
$ws.Range("N7").Value = -2225385

$ws.Range("H8").Value = 7501903
$ws.Range("I8").Value = 14286286
$ws.Range("J8").Value = 2225161
$ws.Range("K8").Value = 14286286
$ws.Range("L8").Value = 2225161
$ws.Range("M8").Value = -14286147
$ws.Range("N8").Value = -2225439

$ws.Range("H122").Value = 175328.14
$ws.Range("I122").Value = 265816.53
$ws.Range("J122").Value = 3400.2
$ws.Range("K122").Value = 797449.5900000001
$ws.Range("L122").Value = 10200.6
$ws.Range("M122").Value = -794999.5900000001
$ws.Range("N122").Value = -15100.6

$ws.Range("H132").Value = 157064.42
$ws.Range("I132").Value = 157064.42
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 471193.26
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -468663.26
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1253.5714
$ws.Range("J2").Value = 1597.5
$ws.Range("L2").Value = 1597.5
$ws.Range("N2").Value = -1821.5

$ws.Range("H7").Value = 6634.5557
$ws.Range("I7").Value = 5078.4443
$ws.Range("J7").Value = 8190.6665
$ws.Range("K7").Value = 5078.4443
$ws.Range("L7").Value = 8190.6665
$ws.Range("M7").Value = -4966.4443
$ws.Range("N7").Value = -8414.666499999999

$ws.Range("H93").Value = 1559.5333
$ws.Range("I93").Value = 1754
$ws.Range("J93").Value = 1024.75
$ws.Range("K93").Value = 1754
$ws.Range("L93").Value = 1024.75
$ws.Range("M93").Value = -506
$ws.Range("N93").Value = -3520.75

$ws.Range("H100").Value = 6541.8213
$ws.Range("I100").Value = 3324.3
$ws.Range("J100").Value = 8329.333000000001
$ws.Range("K100").Value = 3324.3
$ws.Range("L100").Value = 8329.333000000001
$ws.Range("M100").Value = -2783.3
$ws.Range("N100").Value = -9411.333000000001

$ws.Range("H126").Value = 6634.5557
$ws.Range("I126").Value = 5078.4443
$ws.Range("J126").Value = 8190.6665
$ws.Range("K126").Value = 15235.3329
$ws.Range("L126").Value = 24571.9995
$ws.Range("M126").Value = -12765.3329
$ws.Range("N126").Value = -29511.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
